# Add two new rows of columns to the CALLING_HISTORY table (K/L columns):
#   row 8: status / int
#   row 9: description / varchar2
# These mirror the existing "Column"/"Type" table layout used by the other
# entity tables on the sheet (USER, ADMIN, SCHEDULE, RECORD, DENTIST, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (border/style) that's already used by the other data
# rows of this table onto the two new cells before filling in their values.
$ws.Range("K7").Copy($ws.Range("K8"))
$ws.Range("L7").Copy($ws.Range("L8"))
$ws.Range("K7").Copy($ws.Range("K9"))
$ws.Range("L7").Copy($ws.Range("L9"))

$ws.Range("K8").Value = "status"
$ws.Range("L8").Value = "int"
$ws.Range("K9").Value = "description"
$ws.Range("L9").Value = "varchar2"

# Move the active selection to the newly added L9 cell.
[void]$ws.Range("L9").Select()
